$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the trailing spacer rows down by 2 (bottom-up to avoid clobbering) ---
# Row 31 -> Row 33
$ws.Range("B31:J31").Copy($ws.Range("B33"))
# Row 23 -> Row 25
$ws.Range("B23:J23").Copy($ws.Range("B25"))
# Rows 19:20 (old blank spacer rows) -> Rows 21:22
$ws.Range("B19:J20").Copy($ws.Range("B21"))

# Clear the now-stale original rows (23 and 31) that were copied away
$ws.Range("B23:J23").Clear()
$ws.Range("B31:J31").Clear()

# --- Populate new data rows 19 and 20 by cloning the format of row 18 ---
$ws.Range("B18:J18").Copy($ws.Range("B19"))
$ws.Range("B18:J18").Copy($ws.Range("B20"))

# Row 19: Sr# 13
$ws.Range("B19").Value = 13
$ws.Range("C19").Value = "Updates based on comments"
$ws.Range("D19").Value = "04 - 12 - 2019"
$ws.Range("E19").Value = "Ongoing"
$ws.Range("F19").Value = 2

# Row 20: Sr# 14
$ws.Range("B20").Value = 14
$ws.Range("C20").Value = "Lot of feature updates and reports update"
$ws.Range("D20").Value = "07 - 12 - 2019"
$ws.Range("E20").Value = "Ongoing"
$ws.Range("F20").Value = 5

# --- Sheet view tidy-up ---
$ws.Range("D14").Select()
$excel.ActiveWindow.ScrollRow = 1
